$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data values --------------------------------------------------
$ws.Range("C2").Value = 200000
$ws.Range("C7").Value = 8000
$ws.Range("C8:N8").Value = 15000

# --- Extend the used range / formatting out to column O ---------------
# Header row (months) gets left + top alignment
$ws.Range("C1:O1").HorizontalAlignment = -4131
$ws.Range("C1:O1").VerticalAlignment = -4160

# Data block rows 2-5 (labels + monthly columns) -> left alignment
$ws.Range("A2:O5").HorizontalAlignment = -4131

# Data block rows 7-12 (labels + monthly columns) -> left alignment
$ws.Range("A7:O12").HorizontalAlignment = -4131

# Column O on the Expenditure row keeps the plain left-aligned look
$ws.Range("O6").HorizontalAlignment = -4131

# --- Re-merge the Expenditure banner across A6:N6 and center it -------
$ws.Range("A6:N6").Merge()
$ws.Range("A6:N6").HorizontalAlignment = -4108

# --- Selection matches the new merged banner ---------------------------
$ws.Range("A6:N6").Select()
